$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Policy weighting change: Temp Policy_Schools-related weight (column V) switched off ---
# Row 7 holds the category weights; V7 (a weight of 1) becomes 0, and the weight total
# in X7 drops from 13 to 12. This ripples into every LockdownEffectiveness (column X) cell
# below, which is the weighted share: SUMPRODUCT(B7:W7, B{row}:W{row}) / X7.
$ws.Cells.Item(7, 22).Value = 0
$ws.Cells.Item(7, 24).Value = 12

# --- Recalculated LockdownEffectiveness values for existing data rows 24-221 ---
$newX = @(0.08333333333333333, 0.08333333333333333, 0.4166666666666667, 0.4166666666666667, 0.4166666666666667, 0.4166666666666667, 0.4166666666666667, 0.4166666666666667, 0.5833333333249999, 0.5833333333249999, 0.5833333333249999, 0.5833333333249999, 0.9166666666583335, 0.9166666666583335, 0.9166666666583335, 0.9166666666583335, 0.9166666666583335, 0.9166666666583335, 0.9166666666583335, 0.9166666666583335, 0.9166666666583335, 0.9166666666583335, 0.9166666666583335, 0.9166666666583335, 0.9166666666583335, 0.9166666666583335, 0.9166666666583335, 0.9166666666583335, 0.9166666666583335, 0.9166666666583335, 0.9166666666583335, 0.9166666666583335, 0.9166666666583335, 0.9166666666583335, 0.9166666666583335, 0.9166666666583335, 0.9166666666583335, 0.9166666666583335, 0.9166666666583335, 0.9166666666583335, 0.5416666666666666, 0.5416666666666666, 0.5416666666666666, 0.5416666666666666, 0.5416666666666666, 0.5416666666666666, 0.5416666666666666, 0.5416666666666666, 0.5416666666666666, 0.5416666666666666, 0.5416666666666666, 0.5416666666666666, 0.5416666666666666, 0.5416666666666666, 0.3333333333333333, 0.3333333333333333, 0.3333333333333333, 0.3333333333333333, 0.3333333333333333, 0.3333333333333333, 0.3333333333333333, 0.3333333333333333, 0.3333333333333333, 0.3333333333333333, 0.3333333333333333, 0.3333333333333333, 0.3333333333333333, 0.3333333333333333, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25, 0.25)
for ($i = 0; $i -lt $newX.Length; $i++) {
    $ws.Cells.Item(24 + $i, 24).Value = $newX[$i]
}

# --- Append 12 new daily policy rows for 9/30/2020 through 10/11/2020 ---
# Each new row repeats the same per-category flags as the prior day (row 221) and then
# gets its own recalculated LockdownEffectiveness value.
$newDates = @("9/30/2020", "10/1/2020", "10/2/2020", "10/3/2020", "10/4/2020", "10/5/2020", "10/6/2020", "10/7/2020", "10/8/2020", "10/9/2020", "10/10/2020", "10/11/2020")
$rowPattern = @(1, 0, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 1, 1, 0)
for ($i = 0; $i -lt $newDates.Length; $i++) {
    $r = 222 + $i

    # Column A: bold/centered/bordered date label, stored as text (matches existing date rows)
    $ws.Cells.Item($r, 1).Formula = '="' + $newDates[$i] + '"'
    $ws.Cells.Item($r, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4163)
    $ws.Cells.Item($r, 1).Font.Bold = $true
    $ws.Cells.Item($r, 1).HorizontalAlignment = -4108
    $ws.Cells.Item($r, 1).VerticalAlignment = -4160
    $ws.Cells.Item($r, 1).Borders.LineStyle = 1

    # Columns B-W: same per-category flag pattern as row 221
    for ($j = 0; $j -lt $rowPattern.Length; $j++) {
        $ws.Cells.Item($r, 2 + $j).Value = $rowPattern[$j]
    }

    # Column X: recalculated LockdownEffectiveness for the new row
    $ws.Cells.Item($r, 24).Value = 0.25
}
